# Adds %DiffH, %DiffD, %DiffA columns (AF, AG, AH) computed from the
# existing DiffH/DiffD/DiffA (AC/AD/AE) and YtrueH/YtrueD/YtrueA (Z/AA/AB)
# columns, expressed as a percentage: %Diff = Diff / Ytrue * 100

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- Header row (row 1): copy the formatting of the existing last header
# cell (AE1) onto the three new header cells, then set their captions.
$ws.Range("AE1").Copy()
$ws.Range("AF1:AH1").PasteSpecial(-4122)

$ws.Range("AF1").Value = "%DiffH"
$ws.Range("AG1").Value = "%DiffD"
$ws.Range("AH1").Value = "%DiffA"

# --- Data rows: compute the percentage-difference columns with formulas,
# then flatten the results down to plain static values (matching the
# rest of the sheet, which stores computed numbers rather than formulas).
$ws.Range("AF2:AF" + $lastRow).Formula = "=AC2/Z2*100"
$ws.Range("AG2:AG" + $lastRow).Formula = "=AD2/AA2*100"
$ws.Range("AH2:AH" + $lastRow).Formula = "=AE2/AB2*100"

$ws.Range("AF2:AH" + $lastRow).Copy()
$ws.Range("AF2:AH" + $lastRow).PasteSpecial(-4163)

$excel.CutCopyMode = 0
